$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 5).Value = '2026-02-28 07:18:15'
$ws.Cells.Item(3, 5).Value = '2026-02-28 07:18:18'
$ws.Cells.Item(3, 8).Value = '''87%'
$ws.Cells.Item(3, 15).Value = '-0.8 °C'
$ws.Cells.Item(4, 5).Value = '2026-02-28 07:18:20'
$ws.Cells.Item(4, 15).Value = '8.1 °C'
$ws.Cells.Item(5, 5).Value = '2026-02-28 07:18:22'
$ws.Cells.Item(5, 14).Value = '-2.1 °C 6:58 TU'
$ws.Cells.Item(5, 15).Value = '-0.4 °C'
$ws.Cells.Item(6, 5).Value = '2026-02-28 07:18:24'
$ws.Cells.Item(6, 10).Value = '1024.2 hPa'
$ws.Cells.Item(6, 14).Value = '9.8 °C 6:37 TU'
$ws.Cells.Item(6, 15).Value = '10.1 °C'
$ws.Cells.Item(7, 5).Value = '2026-02-28 07:18:26'
$ws.Cells.Item(7, 9).Value = '0.1 mm'
$ws.Cells.Item(7, 15).Value = '12.4 °C'
$ws.Cells.Item(8, 5).Value = '2026-02-28 07:18:28'
$ws.Cells.Item(8, 10).Value = '1023.9 hPa'
$ws.Cells.Item(8, 13).Value = '9.1 °C 6:57 TU'
$ws.Cells.Item(8, 15).Value = '8.6 °C'
$ws.Cells.Item(9, 5).Value = '2026-02-28 07:18:31'
$ws.Cells.Item(9, 14).Value = '5.9 °C 6:34 TU'
$ws.Cells.Item(10, 5).Value = '2026-02-28 07:18:32'
$ws.Cells.Item(10, 13).Value = '10.0 °C 6:59 TU'
$ws.Cells.Item(10, 15).Value = '8.1 °C'
$ws.Cells.Item(11, 5).Value = '2026-02-28 07:18:33'
$ws.Cells.Item(11, 8).Value = '''89%'
$ws.Cells.Item(11, 15).Value = '2.9 °C'
$ws.Cells.Item(12, 5).Value = '2026-02-28 07:18:34'
$ws.Cells.Item(13, 5).Value = '2026-02-28 07:18:35'
$ws.Cells.Item(13, 10).Value = '1026.6 hPa'
$ws.Cells.Item(13, 11).Value = '0.0 MJ/m2'
$ws.Cells.Item(13, 15).Value = '0.8 °C'
$ws.Cells.Item(14, 5).Value = '2026-02-28 07:18:36'
$ws.Cells.Item(14, 8).Value = '''94%'
$ws.Cells.Item(14, 13).Value = '12.8 °C 6:33 TU'
$ws.Cells.Item(14, 15).Value = '10.8 °C'
$ws.Cells.Item(15, 5).Value = '2026-02-28 07:18:37'
$ws.Cells.Item(16, 5).Value = '2026-02-28 07:18:38'
$ws.Cells.Item(16, 11).Value = '0.0 MJ/m2'
$ws.Cells.Item(17, 5).Value = '2026-02-28 07:18:39'
$ws.Cells.Item(17, 8).Value = '''50%'
$ws.Cells.Item(17, 11).Value = '0.0 MJ/m2'
$ws.Cells.Item(17, 15).Value = '4.2 °C'
$ws.Cells.Item(18, 5).Value = '2026-02-28 07:18:40'
$ws.Cells.Item(18, 12).Value = '4.7 km/h - 259º 6:45 TU'
$ws.Cells.Item(18, 13).Value = '9.9 °C 6:54 TU'
$ws.Cells.Item(18, 15).Value = '8.6 °C'
$ws.Cells.Item(19, 5).Value = '2026-02-28 07:18:41'
$ws.Cells.Item(19, 8).Value = '''70%'
$ws.Cells.Item(19, 15).Value = '7.6 °C'
$ws.Cells.Item(20, 5).Value = '2026-02-28 07:18:42'
$ws.Cells.Item(20, 8).Value = '''38%'
$ws.Cells.Item(20, 14).Value = '-1.6 °C 6:42 TU'
$ws.Cells.Item(20, 15).Value = '-0.2 °C'
$ws.Cells.Item(21, 5).Value = '2026-02-28 07:18:44'
$ws.Cells.Item(21, 8).Value = '''78%'
$ws.Cells.Item(21, 11).Value = '0.0 MJ/m2'
$ws.Cells.Item(22, 5).Value = '2026-02-28 07:18:46'
$ws.Cells.Item(22, 14).Value = '-2.1 °C 6:59 TU'
$ws.Cells.Item(23, 5).Value = '2026-02-28 07:18:48'
$ws.Cells.Item(23, 8).Value = '''66%'
$ws.Cells.Item(23, 11).Value = '0.0 MJ/m2'
$ws.Cells.Item(24, 5).Value = '2026-02-28 07:18:51'
$ws.Cells.Item(24, 15).Value = '6.5 °C'
$ws.Cells.Item(25, 5).Value = '2026-02-28 07:18:53'
$ws.Cells.Item(25, 8).Value = '''53%'
$ws.Cells.Item(25, 14).Value = '-0.8 °C 6:30 TU'
$ws.Cells.Item(26, 5).Value = '2026-02-28 07:18:55'
$ws.Cells.Item(26, 10).Value = '1024.1 hPa'
$ws.Cells.Item(26, 15).Value = '4.2 °C'
$ws.Cells.Item(27, 5).Value = '2026-02-28 07:18:58'
$ws.Cells.Item(27, 15).Value = '1.9 °C'
$ws.Cells.Item(28, 5).Value = '2026-02-28 07:19:00'
$ws.Cells.Item(28, 12).Value = '6.1 km/h - 281º 6:57 TU'
$ws.Cells.Item(28, 14).Value = '5.4 °C 6:36 TU'
$ws.Cells.Item(29, 5).Value = '2026-02-28 07:19:02'
$ws.Cells.Item(29, 13).Value = '10.6 °C 6:58 TU'
$ws.Cells.Item(30, 5).Value = '2026-02-28 07:19:04'
$ws.Cells.Item(31, 5).Value = '2026-02-28 07:19:06'
$ws.Cells.Item(31, 8).Value = '''94%'
$ws.Cells.Item(31, 10).Value = '1023.8 hPa'
$ws.Cells.Item(31, 14).Value = '9.5 °C 6:31 TU'
$ws.Cells.Item(32, 5).Value = '2026-02-28 07:19:08'
$ws.Cells.Item(32, 8).Value = '''93%'
$ws.Cells.Item(32, 9).Value = '0.1 mm'
$ws.Cells.Item(33, 5).Value = '2026-02-28 07:19:11'
$ws.Cells.Item(34, 5).Value = '2026-02-28 07:19:13'
$ws.Cells.Item(34, 14).Value = '-2.4 °C 6:57 TU'
$ws.Cells.Item(34, 15).Value = '-0.4 °C'
$ws.Cells.Item(35, 5).Value = '2026-02-28 07:19:15'
$ws.Cells.Item(35, 10).Value = '1023.2 hPa'
$ws.Cells.Item(35, 12).Value = '29.5 km/h - 274º 6:44 TU'
$ws.Cells.Item(35, 14).Value = '4.5 °C 6:48 TU'
$ws.Cells.Item(35, 15).Value = '6.3 °C'
$ws.Cells.Item(36, 5).Value = '2026-02-28 07:19:17'
$ws.Cells.Item(36, 12).Value = '22.0 km/h - 320º 6:43 TU'
$ws.Cells.Item(36, 13).Value = '11.8 °C 6:47 TU'
$ws.Cells.Item(36, 15).Value = '9.9 °C'
$ws.Cells.Item(37, 5).Value = '2026-02-28 07:19:19'
$ws.Cells.Item(37, 14).Value = '2.7 °C 6:36 TU'
$ws.Cells.Item(37, 15).Value = '4.3 °C'
$ws.Cells.Item(38, 5).Value = '2026-02-28 07:19:21'
$ws.Cells.Item(38, 13).Value = '10.0 °C 6:53 TU'
$ws.Cells.Item(39, 5).Value = '2026-02-28 07:19:23'
$ws.Cells.Item(39, 8).Value = '''45%'
$ws.Cells.Item(39, 15).Value = '0.0 °C'
$ws.Cells.Item(40, 5).Value = '2026-02-28 07:19:25'
$ws.Cells.Item(40, 10).Value = '1024.9 hPa'
$ws.Cells.Item(40, 15).Value = '3.1 °C'
$ws.Cells.Item(41, 5).Value = '2026-02-28 07:19:28'
$ws.Cells.Item(41, 8).Value = '''74%'
$ws.Cells.Item(41, 10).Value = '1023.3 hPa'
$ws.Cells.Item(42, 5).Value = '2026-02-28 07:19:30'
$ws.Cells.Item(42, 15).Value = '7.6 °C'
$ws.Cells.Item(43, 5).Value = '2026-02-28 07:19:32'
$ws.Cells.Item(43, 8).Value = '''88%'
$ws.Cells.Item(43, 15).Value = '3.5 °C'
$ws.Cells.Item(44, 5).Value = '2026-02-28 07:19:35'
$ws.Cells.Item(44, 8).Value = '''94%'
$ws.Cells.Item(44, 15).Value = '-1.2 °C'
$ws.Cells.Item(45, 5).Value = '2026-02-28 07:19:37'
$ws.Cells.Item(45, 8).Value = '''91%'
$ws.Cells.Item(45, 10).Value = '1024.4 hPa'
$ws.Cells.Item(45, 14).Value = '4.8 °C 6:59 TU'
$ws.Cells.Item(45, 15).Value = '6.4 °C'
$ws.Cells.Item(46, 5).Value = '2026-02-28 07:19:40'
$ws.Cells.Item(46, 10).Value = '1023.4 hPa'
